$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) holds free-form numeric-looking text (e.g. "29.364.84",
# "0.00001028", "17.60") that must stay literal text rather than being
# auto-coerced to a Number by Excels input parser (which would mangle
# multi-dot values and drop trailing zeros). Pre-format touched D cells as
# Text so the assigned string is preserved exactly.
$ws.Range("D2:D5").NumberFormat = "@"
$ws.Range("D7:D10").NumberFormat = "@"
$ws.Range("D12:D16").NumberFormat = "@"
$ws.Range("D18:D19").NumberFormat = "@"
$ws.Range("D21:D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39:D40").NumberFormat = "@"
$ws.Range("D42:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.364.84"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "1.897.31"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "324.68"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "0.4777"
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("D8").Value = "0.4059"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").Value = "0.08046"
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  -1.34%  "
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "1.947.05"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "5.945"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "7.058"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").Value = "89.66"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("D18").Value = "0.00001028"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "17.60"
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").Value = "29.378.79"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").Value = "5.529"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "11.70"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "2.157"
$ws.Range("E24").Value = "  -1.93%  "
$ws.Range("D25").Value = "2.119.39"
$ws.Range("E25").Value = "  -2.05%  "
$ws.Range("D26").Value = "154.85"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "19.72"
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "6.051"
$ws.Range("E28").Value = "  +5.41%  "
$ws.Range("D29").Value = "2.089"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "118.02"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").Value = "1.020"
$ws.Range("E31").Value = "  -5.14%  "
$ws.Range("D32").Value = "0.09500"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "1.387"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").Value = "3.526"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").Value = "5.378"
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("D37").Value = "0.06051"
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").Value = "0.5851"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "7.840"
$ws.Range("E40").Value = "  -7.10%  "
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "10.13"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").Value = "2.423"
$ws.Range("E43").Value = "  +2.85%  "
$ws.Range("D44").Value = "1.288"
$ws.Range("E44").Value = "  +2.26%  "
$ws.Range("D45").Value = "0.07718"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").Value = "12.18"
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "0.5506"
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("D48").Value = "1.919"
$ws.Range("D49").Value = "112.92"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("D50").Value = "0.2949"
$ws.Range("E50").Value = "  -1.83%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "71.87"
$ws.Range("E51").Value = "  +0.43%  "

Write-Host "Updated cryptos list"
